# Reduce the number of tickers in the 'analise_fundamentalista.xlsx' spreadsheet,
# keeping only the header ("Ticker") plus BBAS3.SA and PETR4.SA, in order to
# reduce the processing time (per commit message).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original rows (1-based):
#  1 Ticker   (header - keep)
#  2 ABCB4.SA (remove)
#  3 ABEV3.SA (remove)
#  4 BBAS3.SA (keep)
#  5 CIEL3.SA (remove)
#  6 CMIG3.SA (remove)
#  7 CMIG4.SA (remove)
#  8 ITSA4.SA (remove)
#  9 ITUB3.SA (remove)
# 10 ODPV3.SA (remove)
# 11 PETR4.SA (keep)
# 12 USIM5.SA (remove)
# 13 WEGE3.SA (remove)
#
# Delete bottom-to-top so row numbers of not-yet-processed rows stay valid.
$rowsToDelete = @(13, 12, 10, 9, 8, 7, 6, 5, 3, 2)
foreach ($r in $rowsToDelete) {
    $ws.Range("A$r").EntireRow.Delete()
}
